# Adição do testng como report e inclusao dos metodos assert
$wb = $excel.ActiveWorkbook

$wsPlanilha1 = $wb.Worksheets.Item("Planilha1")
$wsTeste2    = $wb.Worksheets.Item("Teste2")

# Teste2!A3: "Computador" -> "HP ELITEBOOK FOLIO"
$wsTeste2.Range("A3").Value = "HP ELITEBOOK FOLIO"

# Planilha1!A2: "GabrielaNomuraa" -> "GabbrielaNomura"
$wsPlanilha1.Range("A2").Value = "GabbrielaNomura"
